$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.028960227966309
$ws.Range("B1").Value = 1.815238833427429
$ws.Range("C1").Value = 1.839377164840698
$ws.Range("D1").Value = 1.936708807945251
$ws.Range("E1").Value = 1.357491135597229
